# Auto-generated PowerShell Excel COM-interop script
# Applies numeric cell updates (and a few cell deletions) to rows across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 7500
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 7500
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 22500
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -25246

$ws.Range("H106").Value = 26527896
$ws.Range("I106").Value = 32111862
$ws.Range("K106").Value = 32111862
$ws.Range("M106").Value = -32111231

$ws.Range("H132").Value = 10715.071
$ws.Range("I132").Value = 10319.272
$ws.Range("J132").Value = 12166.333
$ws.Range("K132").Value = 30957.816
$ws.Range("L132").Value = 36498.999
$ws.Range("M132").Value = -28427.816
$ws.Range("N132").Value = -41558.999

$ws.Range("H137").Value = 1327.2727
$ws.Range("I137").Value = 1480
$ws.Range("J137").Value = 1200
$ws.Range("K137").Value = 4440
$ws.Range("L137").Value = 3600
$ws.Range("M137").Value = -1890
$ws.Range("N137").Value = -8700

$ws.Range("H138").Value = 2738.0322
$ws.Range("I138").Value = 2943.1538
$ws.Range("J138").Value = 2683.6123
$ws.Range("K138").Value = 8829.4614
$ws.Range("L138").Value = 8050.836899999999
$ws.Range("M138").Value = -3689.4614
$ws.Range("N138").Value = -18330.8369

$ws.Range("H141").Value = 4581.2856
$ws.Range("I141").Value = 1910.1578
$ws.Range("J141").Value = 10220.333
$ws.Range("K141").Value = 5730.4734
$ws.Range("L141").Value = 30660.999
$ws.Range("M141").Value = -550.4733999999999
$ws.Range("N141").Value = -41020.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 687714.5600000001
$ws.Range("I32").Value = 777287.5600000001
$ws.Range("K32").Value = 777287.5600000001
$ws.Range("M32").Value = -777000.5600000001

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H132").Value = 4604.2
$ws.Range("I132").Value = 5297.7827
$ws.Range("J132").Value = 3274.8333
$ws.Range("K132").Value = 15893.3481
$ws.Range("L132").Value = 9824.499899999999
$ws.Range("M132").Value = -13363.3481
$ws.Range("N132").Value = -14884.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1159.4
$ws.Range("I94").Value = 1132.3334
$ws.Range("K94").Value = 1132.3334
$ws.Range("M94").Value = -681.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 6625
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 6625
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H31").Value = 3578.1924
$ws.Range("I31").Value = 1249.1904
$ws.Range("J31").Value = 5155.9033
$ws.Range("K31").Value = 1249.1904
$ws.Range("L31").Value = 5155.9033
$ws.Range("M31").Value = -954.1904
$ws.Range("N31").Value = -5745.9033

$ws.Range("H34").Value = 3578.1924
$ws.Range("I34").Value = 1249.1904
$ws.Range("J34").Value = 5155.9033
$ws.Range("K34").Value = 1249.1904
$ws.Range("L34").Value = 5155.9033
$ws.Range("M34").Value = -1047.1904
$ws.Range("N34").Value = -5559.9033

$ws.Range("H132").Value = 7248728.5
$ws.Range("I132").Value = 1897.4546
$ws.Range("J132").Value = 13891657
$ws.Range("K132").Value = 5692.3638
$ws.Range("L132").Value = 41674971
$ws.Range("M132").Value = -3162.3638
$ws.Range("N132").Value = -41680031

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 4049.8333
$ws.Range("J41").Value = 4749.75
$ws.Range("L41").Value = 14249.25
$ws.Range("N41").Value = -14925.25

$ws.Range("H68").Value = 1390.9587
$ws.Range("I68").Value = 774.7857
$ws.Range("J68").Value = 1641
$ws.Range("K68").Value = 2324.3571
$ws.Range("L68").Value = 4923
$ws.Range("M68").Value = -1513.3571
$ws.Range("N68").Value = -6545

$ws.Range("H71").Value = 1390.9587
$ws.Range("I71").Value = 774.7857
$ws.Range("J71").Value = 1641
$ws.Range("K71").Value = 6973.071300000001
$ws.Range("L71").Value = 14769
$ws.Range("M71").Value = -2917.071300000001
$ws.Range("N71").Value = -22881

$ws.Range("H104").Value = 3932.7144
$ws.Range("J104").Value = 3932.7144
$ws.Range("L104").Value = 11798.1432
$ws.Range("N104").Value = -17040.1432

$ws.Range("H107").Value = 1503.8513
$ws.Range("I107").Value = 310.54544
$ws.Range("J107").Value = 2464.3171
$ws.Range("K107").Value = 931.63632
$ws.Range("L107").Value = 7392.951300000001
$ws.Range("M107").Value = 988.36368
$ws.Range("N107").Value = -11232.9513

$ws.Range("H131").Value = 1176.6666
$ws.Range("J131").Value = 1187.1428
$ws.Range("L131").Value = 3561.4284
$ws.Range("N131").Value = -13641.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1649.1428
$ws.Range("I97").Value = 1649.1428
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1649.1428
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1153.1428
$ws.Range("N97").ClearContents()

$ws.Range("H132").Value = 2537.75
$ws.Range("I132").Value = 2185.0667
$ws.Range("J132").Value = 3595.8
$ws.Range("K132").Value = 6555.2001
$ws.Range("L132").Value = 10787.4
$ws.Range("M132").Value = -4025.2001
$ws.Range("N132").Value = -15847.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 10290.083
$ws.Range("I93").Value = 15539.143
$ws.Range("J93").Value = 2941.4
$ws.Range("K93").Value = 15539.143
$ws.Range("L93").Value = 2941.4
$ws.Range("M93").Value = -14291.143
$ws.Range("N93").Value = -5437.4

$ws.Range("H132").Value = 3008.1282
$ws.Range("J132").Value = 3060
$ws.Range("L132").Value = 9180
$ws.Range("N132").Value = -14240

$ws.Range("H136").Value = 5377491.5
$ws.Range("I136").Value = 1097.0869
$ws.Range("J136").Value = 20834624
$ws.Range("K136").Value = 3291.2607
$ws.Range("L136").Value = 62503872
$ws.Range("M136").Value = -741.2606999999998
$ws.Range("N136").Value = -62508972

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 37487.25
$ws.Range("I21").Value = 4957.5
$ws.Range("K21").Value = 4957.5
$ws.Range("M21").Value = -4722.5

$ws.Range("H35").Value = 37487.25
$ws.Range("I35").Value = 4957.5
$ws.Range("K35").Value = 4957.5
$ws.Range("M35").Value = -4667.5

$ws.Range("H123").Value = 24730.23
$ws.Range("J123").Value = 24730.23
$ws.Range("L123").Value = 24730.23
$ws.Range("N123").Value = -34530.23

$ws.Range("H132").Value = 4945186.5
$ws.Range("I132").Value = 1539.0322
$ws.Range("J132").Value = 10418510
$ws.Range("K132").Value = 4617.096600000001
$ws.Range("L132").Value = 31255530
$ws.Range("M132").Value = -2087.096600000001
$ws.Range("N132").Value = -31260590

Write-Host "Updated 30 rows across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR."
